# moreinputNY.xlsx - "test-3 (more cities in input, w/o constr."
# Insert "lat" / "lon" columns between the existing data columns and the
# trailing "saldo" column, and populate them per-city.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at K:L. This shifts the existing "saldo" column
# (K) -> M, along with its data/styles, and grows the sheet dimension.
$ws.Columns("K:L").Insert()

# New column headers, styled like the other bold/centered header cells.
$ws.Range("K1").Value2 = "lat"
$ws.Range("L1").Value2 = "lon"
$ws.Range("K1:L1").Font.Bold = $true
$ws.Range("K1:L1").HorizontalAlignment = -4108
$ws.Range("K1:L1").VerticalAlignment = -4160

# Lat/lon values per contiguous city block (rows 2-27).
$latLonByRowRange = @(
    @(2, 2, 47.06, 39.200000000000003),    # Азов
    @(3, 12, 59.8, 33.5),                  # Бокситогорский МР
    @(13, 18, 60.4, 28.4),                 # Выборгский МР
    @(19, 21, 59.3, 30.07),                # Гатчина
    @(22, 26, 59.2, 30.05),                # Гатчинский МР
    @(27, 27, 59.07, 28.05)                # Сланцы
)

foreach ($block in $latLonByRowRange) {
    $startRow = $block[0]
    $endRow = $block[1]
    $lat = $block[2]
    $lon = $block[3]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 11).Value2 = $lat
        $ws.Cells.Item($r, 12).Value2 = $lon
    }
}

# The column insert materialised blank K/L cells on the two trailing
# placeholder rows (32 and 34) that never had data past column J. Clear
# them fully (contents + formatting) so those rows stay untouched, as in
# the authored edit.
$ws.Range("K32:L32").Clear() | Out-Null
$ws.Range("K34:L34").Clear() | Out-Null

# Restore the saved selection recorded in the edited workbook.
$ws.Range("Q13").Select() | Out-Null
